$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.770.19'
$ws.Range("E2").Value = '  +3.78%  '
$ws.Range("D3").Value = '2.730.80'
$ws.Range("E3").Value = '  +2.83%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '581.00'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '157.79'
$ws.Range("E6").Value = '  +8.50%  '
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +4.32%  '
$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '2.755.12'
$ws.Range("E9").Value = '  +3.35%  '
$ws.Range("E10").Value = '  +3.19%  '
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '3.231.25'
$ws.Range("E14").Value = '  +3.33%  '
$ws.Range("D15").Value = '27.37'
$ws.Range("E15").Value = '  +3.08%  '
$ws.Range("D16").Value = '63.727.88'
$ws.Range("E16").Value = '  +3.72%  '
$ws.Range("E17").Value = '  +6.00%  '
$ws.Range("D18").Value = '2.756.96'
$ws.Range("E18").Value = '  +3.13%  '
$ws.Range("D19").Value = '12.09'
$ws.Range("E19").Value = '  +3.93%  '
$ws.Range("D20").Value = '4.95'
$ws.Range("E20").Value = '  +3.67%  '
$ws.Range("D21").Value = '364.31'
$ws.Range("E21").Value = '  +2.49%  '
$ws.Range("D22").Value = '7.02'
$ws.Range("E22").Value = '  +1.37%  '
$ws.Range("D23").Value = '0.541'
$ws.Range("E23").Value = '  +2.73%  '
$ws.Range("D24").Value = '0.997'
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").Value = '66.82'
$ws.Range("E25").Value = '  +4.37%  '
$ws.Range("D26").Value = '0.172'
$ws.Range("E26").Value = '  +5.78%  '
$ws.Range("D27").Value = '8.64'
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").Value = '0.0₃0922'
$ws.Range("E29").Value = '  +11.25%  '
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +6.68%  '
$ws.Range("D32").Value = '1.25'
$ws.Range("E32").Value = '  +13.39%  '
$ws.Range("D33").Value = '173.63'
$ws.Range("E33").Value = '  +4.06%  '
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '20.66'
$ws.Range("E35").Value = '  +2.92%  '
$ws.Range("D36").Value = '4.95'
$ws.Range("E36").Value = '  +4.86%  '
$ws.Range("E37").Value = '  +8.54%  '
$ws.Range("E38").Value = '  +5.61%  '
$ws.Range("E39").Value = '  +10.57%  '
$ws.Range("D40").Value = '4.31'
$ws.Range("E40").Value = '  +4.26%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '6.29'
$ws.Range("E41").Value = '  +16.38%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '340.50'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = '39.48'
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("E44").Value = '  +5.76%  '
$ws.Range("E45").Value = '  +6.64%  '
$ws.Range("D46").Value = '0.0601'
$ws.Range("E46").Value = '  +3.49%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.648'
$ws.Range("E47").Value = '  +3.27%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0260'
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("D49").Value = '138.05'
$ws.Range("E49").Value = '  +1.93%  '
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").Value = '0.995'
$ws.Range("E51").Value = '  +0.07%  '
